# populationData.xlsx update
# Adds three derived columns (Y: PreChangeComb, Z: Years Legal (12 if not legal),
# AA: Change by Year), backfills formulas for V33/V34 (PreLegalChange), and widens
# the two new columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1) -------------------------------------------------
$ws.Range("Y1").Value = "PreChangeComb"
$ws.Range("Z1").Value = "Years Legal (12 if not legal)"
$ws.Range("AA1").Value = "Change by Year"

# --- V33 / V34 gain explicit formulas (previously literal values) ------------
$ws.Range("V33").Formula = "=Q33-F33"
$ws.Range("V34").Formula = "=Q34-F34"

# --- Column Y: PreChangeComb --------------------------------------------------
$ws.Range("Y2").Formula = "=IF(V2="""",W2,V2)"
$ws.Range("Y3:Y52").Formula = "=IF(V3="""",W3,V3)"

# --- Column Z: Years Legal (12 if not legal) ----------------------------------
$ws.Range("Z2").Formula = "=IF(T2<>"""",(2022-T2),12)"
$ws.Range("Z3").Formula = "=IF(T3<>"""",(2022-T3),12)"
$ws.Range("Z4:Z52").Formula = "=IF(T4<>"""",(2022-T4),12)"

# --- Column AA: Change by Year -------------------------------------------------
$ws.Range("AA2").Formula = "=W2/Z2"
$ws.Range("AA3:AA52").Formula = "=W3/Z3"

# --- Column X re-fill so the shared-formula group splits at the row 33/34
#     boundary (matches where V33/V34 picked up explicit formulas). --------
$ws.Range("X2:X33").Formula = "=IF(U2<>"""", U2, W2)"
$ws.Range("X34:X52").Formula = "=IF(U34<>"""", U34, W34)"

# --- Column widths for the two brand-new columns -------------------------------
$ws.Columns.Item(26).ColumnWidth = 17
$ws.Columns.Item(27).ColumnWidth = 15.83

# --- Restore the view state (frozen pane + selection) -------------------------
$ws.Range("AB14").Select()
